# Update cryptos list (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay text even when the replacement string looks
    # like a number (e.g. "632.71"), matching the source workbook's
    # inline-string cells, then drop back to the default style so no
    # stray "Text" number format sticks around on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Column D (Price) updates ---
Set-TextValue "D2"  "69.237.88"
Set-TextValue "D3"  "3.774.95"
Set-TextValue "D5"  "632.71"
Set-TextValue "D6"  "166.92"
Set-TextValue "D7"  "3.771.91"
Set-TextValue "D9"  "0.521"
Set-TextValue "D10" "0.158"
Set-TextValue "D12" "6.73"
Set-TextValue "D13" "0.0000239"
Set-TextValue "D14" "35.12"
Set-TextValue "D15" "4.410.33"
Set-TextValue "D16" "3.756.80"
Set-TextValue "D17" "69.202.47"
Set-TextValue "D18" "17.61"
Set-TextValue "D19" "0.113"
Set-TextValue "D20" "7.02"
Set-TextValue "D21" "463.88"
Set-TextValue "D22" "9.53"
Set-TextValue "D24" "82.71"
Set-TextValue "D25" "0.0000144"
Set-TextValue "D26" "12.09"
Set-TextValue "D28" "10.07"
Set-TextValue "D30" "3.923.98"
Set-TextValue "D31" "2.33"
Set-TextValue "D33" "7.09"
Set-TextValue "D34" "28.48"
Set-TextValue "D35" "0.999"
Set-TextValue "D36" "0.166"
Set-TextValue "D37" "3.726.86"
Set-TextValue "D38" "8.96"
Set-TextValue "D40" "3.29"
Set-TextValue "D50" "46.76"

# --- Column E (Volume(1h)) updates ---
$ws.Range("E2").Value  = "  +1.18%  "
$ws.Range("E3").Value  = "  -0.74%  "
$ws.Range("E4").Value  = "  -0.46%  "
$ws.Range("E5").Value  = "  +3.78%  "
$ws.Range("E6").Value  = "  +2.29%  "
$ws.Range("E7").Value  = "  -0.74%  "
$ws.Range("E8").Value  = "  -0.05%  "
$ws.Range("E9").Value  = "  +0.97%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  +6.05%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +14.23%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  +0.58%  "

# --- Rows 45/46: Monero and Stacks swapped order, with new values ---
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "2.00"
$ws.Range("E45").Value = "  +7.50%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D46" "157.90"
$ws.Range("E46").Value = "  +3.15%  "
